$wb = $excel.ActiveWorkbook

# Set the selection on the "3d" sheet to the full data range (A1:D19)
$ws3d = $wb.Worksheets.Item("3d")
$ws3d.Range("A1:D19").Select() | Out-Null

# Add the new "position" sheet at the end of the workbook
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsPosition = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsPosition.Name = "position"

# Copy the "3d" sheet's data (A1:D19) into the new sheet, offset to start at D3
$srcRows = 19
$srcCols = 4
for ($r = 1; $r -le $srcRows; $r++) {
  for ($c = 1; $c -le $srcCols; $c++) {
    $value = $ws3d.Cells.Item($r, $c).Value2
    $wsPosition.Cells.Item($r + 2, $c + 3).Value = $value
  }
}

# Select J10 on the new sheet, and make it the active sheet/tab
$wsPosition.Range("J10").Select() | Out-Null
